# Update raw data for 2022: append a new change-log entry as row 23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column: match the existing date formatting used by the column
# (numFmtId 15, "d-mmm-yy") before assigning the value so the new cell
# reuses the same style as the rest of the column instead of creating a
# brand-new cell style.
$ws.Range("A23").NumberFormat = "d-mmm-yy"
$ws.Range("A23").Value = (Get-Date -Year 2022 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("B23").Value = 2022
$ws.Range("C23").Value = "Foggy Cove HIGH Quadrat 8"
$ws.Range("D23").Value = "Changed barnacle cover from 440% to 44%"

# Move the selection to follow the newly added row, mirroring the
# author's saved cursor position.
[void]$ws.Range("E23").Select()
